$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 'Bitcoin'
$ws.Range("C2").Value = 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc'
$ws.Range("D2").Value = '38.655.23'
$ws.Range("E2").Value = '  -5.46%  '

$ws.Range("B3").Value = 'Ethereum'
$ws.Range("C3").Value = 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth'
$ws.Range("D3").Value = '2.179.60'
$ws.Range("E3").Value = '  -8.70%  '

$ws.Range("B4").Value = 'TetherUSD'
$ws.Range("C4").Value = 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt'
$ws.Range("D4").Value = '''1.00'
$ws.Range("E4").Value = '  +0.04%  '

$ws.Range("B5").Value = 'BNB'
$ws.Range("C5").Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range("D5").Value = '''291.45'
$ws.Range("E5").Value = '  -7.12%  '

$ws.Range("B6").Value = 'Solana'
$ws.Range("C6").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D6").Value = '''79.22'
$ws.Range("E6").Value = '  -10.17%  '

$ws.Range("B7").Value = 'XRP'
$ws.Range("C7").Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range("D7").Value = '''0.499'
$ws.Range("E7").Value = '  -6.16%  '

$ws.Range("B8").Value = 'USDC'
$ws.Range("C8").Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range("D8").Value = '''1.00'
$ws.Range("E8").Value = '  +0.09%  '

$ws.Range("B9").Value = 'Cardano'
$ws.Range("C9").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D9").Value = '''0.450'
$ws.Range("E9").Value = '  -9.48%  '

$ws.Range("B10").Value = 'Dogecoin'
$ws.Range("C10").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D10").Value = '''0.0763'
$ws.Range("E10").Value = '  -9.22%  '

$ws.Range("B11").Value = 'Avalanche'
$ws.Range("C11").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D11").Value = '''27.43'
$ws.Range("E11").Value = '  -11.81%  '

$ws.Range("B12").Value = 'OKB'
$ws.Range("C12").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D12").Value = '''45.84'
$ws.Range("E12").Value = '  -13.24%  '

$ws.Range("B13").Value = 'TRON'
$ws.Range("C13").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D13").Value = '''0.107'
$ws.Range("E13").Value = '  -2.67%  '

$ws.Range("B14").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C14").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D14").Value = '2.525.46'
$ws.Range("E14").Value = '  -8.39%  '

$ws.Range("B15").Value = 'Polkadot'
$ws.Range("C15").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D15").Value = '''6.01'
$ws.Range("E15").Value = '  -9.05%  '

$ws.Range("B16").Value = 'Chainlink'
$ws.Range("C16").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D16").Value = '''13.62'
$ws.Range("E16").Value = '  -10.59%  '

$ws.Range("B17").Value = 'WrappedEther'
$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D17").Value = '2.208.74'
$ws.Range("E17").Value = '  -7.69%  '

$ws.Range("B18").Value = 'Polygon'
$ws.Range("C18").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D18").Value = '''0.696'
$ws.Range("E18").Value = '  -9.37%  '

$ws.Range("B19").Value = 'WrappedBTC'
$ws.Range("C19").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D19").Value = '38.588.90'
$ws.Range("E19").Value = '  -5.53%  '

$ws.Range("B20").Value = 'ShibaInu'
$ws.Range("C20").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D20").Value = '0.0₃0847'
$ws.Range("E20").Value = '  -7.76%  '

$ws.Range("B21").Value = 'Uniswap'
$ws.Range("C21").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D21").Value = '''5.57'
$ws.Range("E21").Value = '  -9.98%  '

$ws.Range("B22").Value = 'Litecoin'
$ws.Range("C22").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D22").Value = '''63.74'
$ws.Range("E22").Value = '  -7.76%  '

$ws.Range("B23").Value = 'InternetComputer(DFINITY)'
$ws.Range("C23").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D23").Value = '''9.60'
$ws.Range("E23").Value = '  -12.73%  '

$ws.Range("B24").Value = 'BitcoinCash'
$ws.Range("C24").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D24").Value = '''220.63'
$ws.Range("E24").Value = '  -6.00%  '

$ws.Range("B25").Value = 'Dai'
$ws.Range("C25").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D25").Value = '''1.00'
$ws.Range("E25").Value = '  +0.08%  '

$ws.Range("B26").Value = 'PancakeSwap'
$ws.Range("C26").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D26").Value = '''2.34'
$ws.Range("E26").Value = '  -12.38%  '

$ws.Range("B27").Value = 'ImmutableX'
$ws.Range("C27").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D27").Value = '''1.71'
$ws.Range("E27").Value = '  -6.39%  '

$ws.Range("B28").Value = 'EthereumClassic'
$ws.Range("C28").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D28").Value = '''21.88'
$ws.Range("E28").Value = '  -8.77%  '

$ws.Range("B29").Value = 'Toncoin'
$ws.Range("C29").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D29").Value = '''2.15'
$ws.Range("E29").Value = '  -2.98%  '

$ws.Range("B30").Value = 'Cosmos'
$ws.Range("C30").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D30").Value = '''8.75'
$ws.Range("E30").Value = '  -7.11%  '

$ws.Range("B31").Value = 'Monero'
$ws.Range("C31").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D31").Value = '''147.18'
$ws.Range("E31").Value = '  -4.05%  '

$ws.Range("B32").Value = 'FirstDigitalUSD'
$ws.Range("C32").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D32").Value = '''1.00'
$ws.Range("E32").Value = '  -0.10%  '

$ws.Range("B33").Value = 'InjectiveProtocol'
$ws.Range("C33").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D33").Value = '''30.52'
$ws.Range("E33").Value = '  -10.39%  '

$ws.Range("B34").Value = 'WEMIXToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D34").Value = '''2.32'
$ws.Range("E34").Value = '  -5.16%  '

$ws.Range("B35").Value = 'Filecoin'
$ws.Range("C35").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D35").Value = '''4.64'
$ws.Range("E35").Value = '  -12.14%  '

$ws.Range("B36").Value = 'Hedera'
$ws.Range("C36").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D36").Value = '''0.0674'
$ws.Range("E36").Value = '  -8.32%  '

$ws.Range("B37").Value = 'Stellar'
$ws.Range("C37").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D37").Value = '''0.107'
$ws.Range("E37").Value = '  -5.80%  '

$ws.Range("B38").Value = 'Kaspa'
$ws.Range("C38").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D38").Value = '''0.0937'
$ws.Range("E38").Value = '  -6.54%  '

$ws.Range("B39").Value = 'LidoDAOToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D39").Value = '''2.57'
$ws.Range("E39").Value = '  -8.95%  '

$ws.Range("B40").Value = 'ARBITRUM'
$ws.Range("C40").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D40").Value = '''1.55'
$ws.Range("E40").Value = '  -10.47%  '

$ws.Range("B41").Value = 'Celestia'
$ws.Range("C41").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D41").Value = '''14.03'
$ws.Range("E41").Value = '  -13.60%  '

$ws.Range("B42").Value = 'RenderToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D42").Value = '''3.52'
$ws.Range("E42").Value = '  -8.78%  '

$ws.Range("B43").Value = 'Maker'
$ws.Range("C43").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D43").Value = '1.874.48'
$ws.Range("E43").Value = '  -4.91%  '

$ws.Range("B44").Value = 'ApeXProtocol'
$ws.Range("C44").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D44").Value = '''2.06'
$ws.Range("E44").Value = '  -12.73%  '

$ws.Range("B45").Value = 'VeChain'
$ws.Range("C45").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D45").Value = '''0.0251'
$ws.Range("E45").Value = '  -7.79%  '

$ws.Range("B46").Value = 'EnergySwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D46").Value = '''15.79'
$ws.Range("E46").Value = '  -12.08%  '

$ws.Range("B47").Value = 'FraxShare'
$ws.Range("C47").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D47").Value = '''8.69'
$ws.Range("E47").Value = '  -10.09%  '

$ws.Range("B48").Value = 'NEARProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D48").Value = '''2.47'
$ws.Range("E48").Value = '  -10.50%  '

$ws.Range("B49").Value = 'RocketPoolETH'
$ws.Range("C49").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D49").Value = '2.415.33'
$ws.Range("E49").Value = '  -7.83%  '

$ws.Range("B50").Value = 'BitcoinSV'
$ws.Range("C50").Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range("D50").Value = '''68.38'
$ws.Range("E50").Value = '  -6.21%  '

$ws.Range("B51").Value = 'Aave'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D51").Value = '''84.92'
$ws.Range("E51").Value = '  -9.77%  '
